$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ganesh PK (row 9) and Amarnath (row 4) paid their July,18 subscription (column G)
$ws.Range("G4").Value = 500
$ws.Range("G9").Value = 500

# Update the selected/active cell to reflect the last edited cell as per the diff
$ws.Range("G4").Select()
